$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the frozen-pane view so the bottom-right pane's top-left cell
# becomes B116 (was B106). Best-effort: some hosts do not persist
# window/pane scroll position across save, but this mirrors what Excel
# does when a user scrolls a frozen sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 116
$excel.ActiveWindow.ScrollColumn = 2

# Updated filesize measurements (columns C and E) for rows 9-120.
$ws.Range("C9").Value = 3295971
$ws.Range("E9").Value = 3295963
$ws.Range("C10").Value = 879397
$ws.Range("E10").Value = 879397
$ws.Range("C11").Value = 501531
$ws.Range("E11").Value = 501529
$ws.Range("C15").Value = 180784
$ws.Range("E15").Value = 180784
$ws.Range("C16").Value = 176001
$ws.Range("E16").Value = 176001
$ws.Range("C18").Value = 3756322
$ws.Range("E18").Value = 3756322
$ws.Range("C19").Value = 433281
$ws.Range("E19").Value = 433281
$ws.Range("C20").Value = 958686
$ws.Range("E20").Value = 958630
$ws.Range("C21").Value = 96557
$ws.Range("E21").Value = 96557
$ws.Range("C22").Value = 2634900
$ws.Range("E22").Value = 2634900
$ws.Range("C23").Value = 232552
$ws.Range("E23").Value = 232552
$ws.Range("C24").Value = 36826
$ws.Range("E24").Value = 36826
$ws.Range("C25").Value = 111204
$ws.Range("E25").Value = 111204
$ws.Range("C26").Value = 192545
$ws.Range("E26").Value = 192545
$ws.Range("C27").Value = 142114
$ws.Range("E27").Value = 142114
$ws.Range("C28").Value = 76519
$ws.Range("E28").Value = 76519
$ws.Range("C29").Value = 239123
$ws.Range("E29").Value = 239123
$ws.Range("C30").Value = 510267
$ws.Range("E30").Value = 510267
$ws.Range("C32").Value = 1712403
$ws.Range("E32").Value = 1712403
$ws.Range("C35").Value = 441503
$ws.Range("E35").Value = 441503
$ws.Range("C36").Value = 188377
$ws.Range("E36").Value = 188377
$ws.Range("C37").Value = 252496
$ws.Range("E37").Value = 252495
$ws.Range("C38").Value = 859214
$ws.Range("E38").Value = 859214
$ws.Range("C40").Value = 270278
$ws.Range("E40").Value = 270278
$ws.Range("C43").Value = 16430
$ws.Range("E43").Value = 16430
$ws.Range("C45").Value = 23156
$ws.Range("E45").Value = 23156
$ws.Range("C49").Value = 197743
$ws.Range("E49").Value = 197743
$ws.Range("C50").Value = 72044
$ws.Range("E50").Value = 72044
$ws.Range("C51").Value = 49374
$ws.Range("E51").Value = 49374
$ws.Range("C55").Value = 13042
$ws.Range("E55").Value = 13042
$ws.Range("C58").Value = 234331
$ws.Range("E58").Value = 234331
$ws.Range("C59").Value = 36935
$ws.Range("E59").Value = 36935
$ws.Range("C60").Value = 85599
$ws.Range("E60").Value = 85552
$ws.Range("C61").Value = 9705
$ws.Range("E61").Value = 9705
$ws.Range("C62").Value = 186523
$ws.Range("E62").Value = 186523
$ws.Range("C63").Value = 17989
$ws.Range("E63").Value = 17989
$ws.Range("C65").Value = 24248
$ws.Range("E65").Value = 24248
$ws.Range("C66").Value = 17371
$ws.Range("E66").Value = 17371
$ws.Range("C67").Value = 20877
$ws.Range("E67").Value = 20877
$ws.Range("C68").Value = 13639
$ws.Range("E68").Value = 13639
$ws.Range("C69").Value = 34863
$ws.Range("E69").Value = 34863
$ws.Range("C70").Value = 97299
$ws.Range("E70").Value = 97299
$ws.Range("C71").Value = 25675
$ws.Range("E71").Value = 25675
$ws.Range("C72").Value = 151112
$ws.Range("E72").Value = 151112
$ws.Range("C74").Value = 30087
$ws.Range("E74").Value = 30087
$ws.Range("C75").Value = 41313
$ws.Range("E75").Value = 41313
$ws.Range("C76").Value = 27222
$ws.Range("E76").Value = 27222
$ws.Range("C77").Value = 45456
$ws.Range("E77").Value = 45456
$ws.Range("C78").Value = 114794
$ws.Range("E78").Value = 114794
$ws.Range("C79").Value = 17198
$ws.Range("E79").Value = 17198
$ws.Range("C80").Value = 22795
$ws.Range("E80").Value = 22795
$ws.Range("C81").Value = 26806
$ws.Range("E81").Value = 26806
$ws.Range("C85").Value = 77903
$ws.Range("E85").Value = 77903
$ws.Range("C89").Value = 839589
$ws.Range("E89").Value = 839589
$ws.Range("C90").Value = 238633
$ws.Range("E90").Value = 238633
$ws.Range("C91").Value = 156224
$ws.Range("E91").Value = 156224
$ws.Range("C95").Value = 47754
$ws.Range("E95").Value = 47754
$ws.Range("C96").Value = 54174
$ws.Range("E96").Value = 54174
$ws.Range("C98").Value = 987077
$ws.Range("E98").Value = 987077
$ws.Range("C99").Value = 116675
$ws.Range("E99").Value = 116675
$ws.Range("C100").Value = 281042
$ws.Range("E100").Value = 281042
$ws.Range("C101").Value = 29220
$ws.Range("E101").Value = 29220
$ws.Range("C102").Value = 703854
$ws.Range("E102").Value = 703854
$ws.Range("C103").Value = 62558
$ws.Range("E103").Value = 62558
$ws.Range("C104").Value = 13697
$ws.Range("E104").Value = 13697
$ws.Range("C105").Value = 50453
$ws.Range("E105").Value = 50453
$ws.Range("C106").Value = 53953
$ws.Range("E106").Value = 53953
$ws.Range("C107").Value = 54793
$ws.Range("E107").Value = 54793
$ws.Range("C108").Value = 32259
$ws.Range("E108").Value = 32259
$ws.Range("C109").Value = 91118
$ws.Range("E109").Value = 91118
$ws.Range("C110").Value = 251633
$ws.Range("E110").Value = 251633
$ws.Range("C111").Value = 75896
$ws.Range("E111").Value = 75896
$ws.Range("C112").Value = 491499
$ws.Range("E112").Value = 491499
$ws.Range("C114").Value = 104837
$ws.Range("E114").Value = 104837
$ws.Range("C115").Value = 135459
$ws.Range("E115").Value = 135459
$ws.Range("C116").Value = 70049
$ws.Range("E116").Value = 70049
$ws.Range("C117").Value = 107848
$ws.Range("E117").Value = 107848
$ws.Range("C118").Value = 313237
$ws.Range("E118").Value = 313237
$ws.Range("C119").Value = 62962
$ws.Range("E119").Value = 62962
$ws.Range("C120").Value = 77200
$ws.Range("E120").Value = 77200
